# edit.ps1 - applies the commit "feat: add 2022-Q1 data"
# 1) Repurpose the current total sheet ("总计") into the new quarter sheet "2022-Q1"
#    (content replaced with the 2022-Q1 fund holdings table), keeping its sheetId/rId.
# 2) Insert a brand new "总计" sheet right after it, containing the historical totals
#    table with a new first row for 2022-Q1 prepended and old rows shifted down.

$wb = $excel.ActiveWorkbook

$quarterTemplate = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# --- Build the new "总计" sheet first (while $oldTotal still holds the old data) ---
$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$newTotal.Name = "TEMP_NEW_TOTAL"

# header row, copied verbatim (keeps text + style)
$oldTotal.Range("A1:D1").Copy($newTotal.Range("A1"))

# row 2 = new 2022-Q1 entry; copy a data-row template for correct styling, then overwrite values
$oldTotal.Range("A2:D2").Copy($newTotal.Range("A2"))
$newTotal.Cells.Item(2, 1).Value = 0
$newTotal.Cells.Item(2, 2).Value = "2022-Q1"
$newTotal.Cells.Item(2, 3).Value = 18
$newTotal.Cells.Item(2, 4).Value = 2.69

# rows 3-6 = old rows 2-5 (2021-Q4 .. 2021-Q1), copied with formatting, then fix the index column
$oldTotal.Range("A2:D5").Copy($newTotal.Range("A3"))
$newTotal.Cells.Item(3, 1).Value = 1
$newTotal.Cells.Item(4, 1).Value = 2
$newTotal.Cells.Item(5, 1).Value = 3
$newTotal.Cells.Item(6, 1).Value = 4

# --- Repurpose the old "总计" sheet into "2022-Q1" ---
# (rename away from "总计" first so the new sheet can take that name without a collision)
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $oldTotal

# Clear the old 总计 content (A1:D5) before laying out the new H-column table
$q1.Range("A1:D5").Clear()

# header row (same column headers used on every other quarter sheet)
$quarterTemplate.Range("A1:H1").Copy($q1.Range("A1"))

# lay down 18 data rows of correctly-styled template (6-row blocks, 3x)
$quarterTemplate.Range("A2:H7").Copy($q1.Range("A2"))
$quarterTemplate.Range("A2:H7").Copy($q1.Range("A8"))
$quarterTemplate.Range("A2:H7").Copy($q1.Range("A14"))

# text-typed columns (fund code / name / scale / position / pct / value) must stay text
$q1.Range("B2:G19").NumberFormat = "@"

$fundData = @(
    @(0, "166301", "华商新趋势优选灵活配置混合", "26.96", "86.39", "2.43", "0.6551", 8),
    @(1, "630002", "华商盛世成长混合", "22.62", "93.39", "2.86", "0.6469", 8),
    @(2, "000390", "华商优势行业混合", "19.61", "88.06", "2.94", "0.5765", 5),
    @(3, "002601", "中银证券价值精选灵活配置混合", "3.41", "93.74", "6.06", "0.2066", 2),
    @(4, "014277", "万家北交所慧选两年定期开放混合A", "4.47", "97.90", "3.80", "0.1699", 6),
    @(5, "005313", "万家中证1000指数增强A", "9.01", "93.72", "1.11", "0.1000", 5),
    @(6, "519183", "万家双引擎灵活配置混合", "2.23", "94.01", "4.42", "0.0986", 10),
    @(7, "005314", "万家中证1000指数增强C", "4.95", "93.72", "1.11", "0.0549", 5),
    @(8, "400007", "东方策略成长混合", "1.85", "82.01", "2.87", "0.0531", 9),
    @(9, "001219", "上投摩根动态多因子策略混合", "1.16", "92.44", "3.28", "0.0380", 4),
    @(10, "310368", "申万菱信竞争优势混合", "0.83", "91.22", "3.47", "0.0288", 10),
    @(11, "001244", "华泰柏瑞量化智慧灵活配置混合A", "3.50", "91.02", "0.76", "0.0266", 10),
    @(12, "014278", "万家北交所慧选两年定期开放混合C", "0.55", "97.90", "3.80", "0.0209", 6),
    @(13, "006104", "华泰柏瑞量化智慧灵活配置混合C", "0.84", "91.02", "0.76", "0.0064", 10),
    @(14, "004794", "富荣福鑫灵活配置混合A", "0.06", "89.60", "3.79", "0.0023", 6),
    @(15, "004795", "富荣福鑫灵活配置混合C", "0.06", "89.60", "3.79", "0.0023", 6),
    @(16, "006857", "蜂巢卓睿灵活配置混合A", "0.14", "78.16", "1.41", "0.0020", 3),
    @(17, "006858", "蜂巢卓睿灵活配置混合C", "0.04", "78.16", "1.41", "0.0006", 3)
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $fundData[$i]
    $row = 2 + $i
    $q1.Cells.Item($row, 1).Value = $r[0]
    $q1.Cells.Item($row, 2).Value = $r[1]
    $q1.Cells.Item($row, 3).Value = $r[2]
    $q1.Cells.Item($row, 4).Value = $r[3]
    $q1.Cells.Item($row, 5).Value = $r[4]
    $q1.Cells.Item($row, 6).Value = $r[5]
    $q1.Cells.Item($row, 7).Value = $r[6]
    $q1.Cells.Item($row, 8).Value = $r[7]
}

$q1.Range("A1").Select()
